# Update column F ("dSF") values on Sheet1 to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -4
    3  = -7
    5  = -1
    9  = 6
    11 = -2
    12 = 4
    13 = -1
    14 = -3
    15 = -2
    22 = -3
    24 = 0
    27 = 3
    28 = -4
    32 = -3
    41 = 6
    42 = -1
    45 = -2
    47 = 5
    51 = 7
    52 = -2
    54 = 4
    57 = -2
    59 = -3
    64 = -7
    67 = 3
    68 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
